# Project Blueprint Template (Games) - Kapitel 2 edit
#
# 1. Replace the "Target Audience" section placeholder paragraph + the
#    German bullet list with the new English "two audience groups" text.
# 2. Move the _GoBack bookmark from the "Competitive Gameplay" heading to
#    the very start of the new Target Audience text.

$d = $word.ActiveDocument

# --- remove the old _GoBack bookmark from "Competitive Gameplay" ---------
# Do this first, while it is still the only bookmark of that name in the
# document, to avoid any ambiguity about which same-named bookmark gets
# targeted once the new one (added further below) exists too.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- locate the block to replace -----------------------------------------
# Find the intro paragraph "Explain who your audience is. ..."
$introRange = $d.Content
$found = $introRange.Find.Execute("Explain who your audience is", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$introRange.Expand(4)  # wdParagraph
$introStart = $introRange.Start

# Find the last bullet paragraph "Perfektionisten" (end of the list block)
$endRange = $d.Content
$found2 = $endRange.Find.Execute("Perfektionisten", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endRange.Expand(4)  # wdParagraph
$blockEnd = $endRange.End

# Delete the whole placeholder block (intro paragraph + every bullet item),
# leaving the trailing empty paragraph (the one right before "Personas")
# untouched.
$deleteRange = $d.Range($introStart, $blockEnd)
$deleteRange.Delete()

# --- insert the new content -----------------------------------------------
$insertPoint = $d.Range($introStart, $introStart)

$newContentXml = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
<w:r><w:t xml:space="preserve">Our Audience is separated into two groups. The first group </w:t></w:r>
<w:r><w:t xml:space="preserve">is made up of students from high schools or universities. </w:t></w:r>
<w:r><w:t xml:space="preserve">This group will be our main Audience. </w:t></w:r>
<w:r><w:t xml:space="preserve">Teachers can use </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>UnrealCup</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> as part of their lessons for example as introduction to show how simple logic affects the behavior of the football players.</w:t></w:r>
<w:r><w:t xml:space="preserve"> The students can also compete against each other and therefore will be encouraged to improve their skills even further. For </w:t></w:r>
<w:r><w:t>an</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t>example of this audience</w:t></w:r>
<w:r><w:t xml:space="preserve"> you can look at the persona &#8220;John </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>Smartman</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t>&#8221;.</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t xml:space="preserve">The second target group are more of a gamer who likes to play strategy games. </w:t></w:r>
<w:r><w:t xml:space="preserve">This group won&#8217;t be as big as the earlier mentioned group. But this group is made out of people who want to create the perfect team, which competes only with the best. This group features </w:t></w:r>
<w:r><w:t xml:space="preserve">people roughly in their mid-20s who are willing to spend a lot of time to create the best team they can. They won&#8217;t stop by using the graphical editor to make their team, but they will use all possibilities available, even if they have to learn something new. An example of this group is &#8220;Jack Hardcore&#8221;, who is mentioned in the personas chapter. </w:t></w:r>
</w:p>
<w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$insertPoint.InsertXML($newContentXml)
